$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 135
$ws.Cells.Item(135, 1).Value = 133
$ws.Cells.Item(135, 2).Value = 7127394
$ws.Cells.Item(135, 3).Value = "Australia ALeague"
$ws.Cells.Item(135, 4).Value = "Australia ALeague"
$ws.Cells.Item(135, 5).Value = 45381.875
$ws.Cells.Item(135, 6).Value = "Wellington Phoenix"
$ws.Cells.Item(135, 7).Value = "Brisbane Roar"
$ws.Cells.Item(135, 8).Value = 1
$ws.Cells.Item(135, 9).Value = 0
$ws.Cells.Item(135, 10).Value = "H"
$ws.Cells.Item(135, 11).Value = 1.8
$ws.Cells.Item(135, 12).Value = 3.8
$ws.Cells.Item(135, 13).Value = 4
$ws.Cells.Item(135, 14).Value = 2.2
$ws.Cells.Item(135, 15).Value = 3.5
$ws.Cells.Item(135, 16).Value = 3.2
$ws.Cells.Item(135, 17).Value = -0.25
$ws.Cells.Item(135, 18).Value = 1.925
$ws.Cells.Item(135, 19).Value = 1.925
$ws.Cells.Item(135, 20).Value = 2.75
$ws.Cells.Item(135, 21).Value = 1.825
$ws.Cells.Item(135, 22).Value = 2.025
$ws.Cells.Item(135, 23).Value = 1.2
$ws.Cells.Item(135, 24).Value = -1
$ws.Cells.Item(135, 25).Value = -1
$ws.Cells.Item(135, 26).Value = 0.925
$ws.Cells.Item(135, 27).Value = -1
$ws.Cells.Item(135, 28).Value = -1
$ws.Cells.Item(135, 29).Value = 1.025

# Row 136
$ws.Cells.Item(136, 1).Value = 134
$ws.Cells.Item(136, 2).Value = 7127397
$ws.Cells.Item(136, 3).Value = "Australia ALeague"
$ws.Cells.Item(136, 4).Value = "Australia ALeague"
$ws.Cells.Item(136, 5).Value = 45382.04166666666
$ws.Cells.Item(136, 6).Value = "Melbourne Victory"
$ws.Cells.Item(136, 7).Value = "Perth Glory"
$ws.Cells.Item(136, 8).Value = 2
$ws.Cells.Item(136, 9).Value = 1
$ws.Cells.Item(136, 10).Value = "H"
$ws.Cells.Item(136, 11).Value = 1.4
$ws.Cells.Item(136, 12).Value = 5
$ws.Cells.Item(136, 13).Value = 6.5
$ws.Cells.Item(136, 14).Value = 1.363
$ws.Cells.Item(136, 15).Value = 6
$ws.Cells.Item(136, 16).Value = 7
$ws.Cells.Item(136, 17).Value = -1.5
$ws.Cells.Item(136, 18).Value = 1.875
$ws.Cells.Item(136, 19).Value = 1.975
$ws.Cells.Item(136, 20).Value = 3.5
$ws.Cells.Item(136, 21).Value = 1.925
$ws.Cells.Item(136, 22).Value = 1.925
$ws.Cells.Item(136, 23).Value = 0.363
$ws.Cells.Item(136, 24).Value = -1
$ws.Cells.Item(136, 25).Value = -1
$ws.Cells.Item(136, 26).Value = -1
$ws.Cells.Item(136, 27).Value = 0.9750000000000001
$ws.Cells.Item(136, 28).Value = -1
$ws.Cells.Item(136, 29).Value = 0.925

# Row 137
$ws.Cells.Item(137, 1).Value = 135
$ws.Cells.Item(137, 2).Value = 7127398
$ws.Cells.Item(137, 3).Value = "Australia ALeague"
$ws.Cells.Item(137, 4).Value = "Australia ALeague"
$ws.Cells.Item(137, 5).Value = 45383.04166666666
$ws.Cells.Item(137, 6).Value = "Macarthur FC"
$ws.Cells.Item(137, 7).Value = "Western Sydney Wanderers"
$ws.Cells.Item(137, 8).Value = 1
$ws.Cells.Item(137, 9).Value = 3
$ws.Cells.Item(137, 10).Value = "A"
$ws.Cells.Item(137, 11).Value = 2.5
$ws.Cells.Item(137, 12).Value = 3.5
$ws.Cells.Item(137, 13).Value = 2.625
$ws.Cells.Item(137, 14).Value = 2.6
$ws.Cells.Item(137, 15).Value = 3.8
$ws.Cells.Item(137, 16).Value = 2.45
$ws.Cells.Item(137, 17).Value = 0
$ws.Cells.Item(137, 18).Value = 1.98
$ws.Cells.Item(137, 19).Value = 1.92
$ws.Cells.Item(137, 20).Value = 3.25
$ws.Cells.Item(137, 21).Value = 1.9
$ws.Cells.Item(137, 22).Value = 1.95
$ws.Cells.Item(137, 23).Value = -1
$ws.Cells.Item(137, 24).Value = -1
$ws.Cells.Item(137, 25).Value = 1.45
$ws.Cells.Item(137, 26).Value = -1
$ws.Cells.Item(137, 27).Value = 0.9199999999999999
$ws.Cells.Item(137, 28).Value = 0.8999999999999999
$ws.Cells.Item(137, 29).Value = -1

# Row 138
$ws.Range("A134").Copy()
$ws.Range("A138").PasteSpecial(-4122)
$ws.Range("E134").Copy()
$ws.Range("E138").PasteSpecial(-4122)
$ws.Cells.Item(138, 1).Value = 136
$ws.Cells.Item(138, 2).Value = 7898681
$ws.Cells.Item(138, 3).Value = "Australia ALeague"
$ws.Cells.Item(138, 4).Value = "Australia ALeague"
$ws.Cells.Item(138, 5).Value = 45384.20833333334
$ws.Cells.Item(138, 6).Value = "Central Coast Mariners"
$ws.Cells.Item(138, 7).Value = "Melbourne City"
$ws.Cells.Item(138, 11).Value = 2.1
$ws.Cells.Item(138, 12).Value = 4
$ws.Cells.Item(138, 13).Value = 3
$ws.Cells.Item(138, 14).Value = 2.2
$ws.Cells.Item(138, 15).Value = 3.8
$ws.Cells.Item(138, 16).Value = 3
$ws.Cells.Item(138, 17).Value = -0.25
$ws.Cells.Item(138, 18).Value = 1.95
$ws.Cells.Item(138, 19).Value = 1.95
$ws.Cells.Item(138, 20).Value = 3
$ws.Cells.Item(138, 21).Value = 2
$ws.Cells.Item(138, 22).Value = 1.85
$ws.Cells.Item(138, 23).Value = 0
$ws.Cells.Item(138, 24).Value = 0
$ws.Cells.Item(138, 25).Value = 0
$ws.Cells.Item(138, 26).Value = 0
$ws.Cells.Item(138, 27).Value = 0

# Row 139
$ws.Range("A134").Copy()
$ws.Range("A139").PasteSpecial(-4122)
$ws.Range("E134").Copy()
$ws.Range("E139").PasteSpecial(-4122)
$ws.Cells.Item(139, 1).Value = 137
$ws.Cells.Item(139, 2).Value = 7661947
$ws.Cells.Item(139, 3).Value = "Australia ALeague"
$ws.Cells.Item(139, 4).Value = "Australia ALeague"
$ws.Cells.Item(139, 5).Value = 45385.32291666666
$ws.Cells.Item(139, 6).Value = "Perth Glory"
$ws.Cells.Item(139, 7).Value = "Sydney FC"
$ws.Cells.Item(139, 11).Value = 3.1
$ws.Cells.Item(139, 12).Value = 3.6
$ws.Cells.Item(139, 13).Value = 2.2
$ws.Cells.Item(139, 14).Value = 3.8
$ws.Cells.Item(139, 15).Value = 4.2
$ws.Cells.Item(139, 16).Value = 1.8
$ws.Cells.Item(139, 17).Value = 0.75
$ws.Cells.Item(139, 18).Value = 1.84
$ws.Cells.Item(139, 19).Value = 2.06
$ws.Cells.Item(139, 20).Value = 3.25
$ws.Cells.Item(139, 21).Value = 1.875
$ws.Cells.Item(139, 22).Value = 1.975
$ws.Cells.Item(139, 23).Value = 0
$ws.Cells.Item(139, 24).Value = 0
$ws.Cells.Item(139, 25).Value = 0
$ws.Cells.Item(139, 26).Value = 0
$ws.Cells.Item(139, 27).Value = 0

# Row 140
$ws.Range("A134").Copy()
$ws.Range("A140").PasteSpecial(-4122)
$ws.Range("E134").Copy()
$ws.Range("E140").PasteSpecial(-4122)
$ws.Cells.Item(140, 1).Value = 138
$ws.Cells.Item(140, 2).Value = 7127399
$ws.Cells.Item(140, 3).Value = "Australia ALeague"
$ws.Cells.Item(140, 4).Value = "Australia ALeague"
$ws.Cells.Item(140, 5).Value = 45387.23958333334
$ws.Cells.Item(140, 6).Value = "Western Sydney Wanderers"
$ws.Cells.Item(140, 7).Value = "Brisbane Roar"
$ws.Cells.Item(140, 11).Value = 2.1
$ws.Cells.Item(140, 12).Value = 3.75
$ws.Cells.Item(140, 13).Value = 3.1
$ws.Cells.Item(140, 14).Value = 2
$ws.Cells.Item(140, 15).Value = 4
$ws.Cells.Item(140, 16).Value = 3.2
$ws.Cells.Item(140, 17).Value = -0.5
$ws.Cells.Item(140, 18).Value = 2.03
$ws.Cells.Item(140, 19).Value = 1.87
$ws.Cells.Item(140, 20).Value = 3.25
$ws.Cells.Item(140, 21).Value = 2.025
$ws.Cells.Item(140, 22).Value = 1.825
$ws.Cells.Item(140, 23).Value = 0
$ws.Cells.Item(140, 24).Value = 0
$ws.Cells.Item(140, 25).Value = 0
$ws.Cells.Item(140, 26).Value = 0
$ws.Cells.Item(140, 27).Value = 0

# Row 141
$ws.Range("A134").Copy()
$ws.Range("A141").PasteSpecial(-4122)
$ws.Range("E134").Copy()
$ws.Range("E141").PasteSpecial(-4122)
$ws.Cells.Item(141, 1).Value = 139
$ws.Cells.Item(141, 2).Value = 8034339
$ws.Cells.Item(141, 3).Value = "Australia ALeague"
$ws.Cells.Item(141, 4).Value = "Australia ALeague"
$ws.Cells.Item(141, 5).Value = 45388.0625
$ws.Cells.Item(141, 6).Value = "Western United FC"
$ws.Cells.Item(141, 7).Value = "Macarthur FC"
$ws.Cells.Item(141, 11).Value = 2.6
$ws.Cells.Item(141, 12).Value = 3.5
$ws.Cells.Item(141, 13).Value = 2.55
$ws.Cells.Item(141, 14).Value = 2.1
$ws.Cells.Item(141, 15).Value = 3.75
$ws.Cells.Item(141, 16).Value = 3.1
$ws.Cells.Item(141, 17).Value = -0.25
$ws.Cells.Item(141, 18).Value = 1.85
$ws.Cells.Item(141, 19).Value = 2.05
$ws.Cells.Item(141, 20).Value = 3.25
$ws.Cells.Item(141, 21).Value = 1.925
$ws.Cells.Item(141, 22).Value = 1.925
$ws.Cells.Item(141, 23).Value = 0
$ws.Cells.Item(141, 24).Value = 0
$ws.Cells.Item(141, 25).Value = 0
$ws.Cells.Item(141, 26).Value = 0
$ws.Cells.Item(141, 27).Value = 0

# Row 142
$ws.Range("A134").Copy()
$ws.Range("A142").PasteSpecial(-4122)
$ws.Range("E134").Copy()
$ws.Range("E142").PasteSpecial(-4122)
$ws.Cells.Item(142, 1).Value = 140
$ws.Cells.Item(142, 2).Value = 8005739
$ws.Cells.Item(142, 3).Value = "Australia ALeague"
$ws.Cells.Item(142, 4).Value = "Australia ALeague"
$ws.Cells.Item(142, 5).Value = 45388.14583333334
$ws.Cells.Item(142, 6).Value = "Central Coast Mariners"
$ws.Cells.Item(142, 7).Value = "Wellington Phoenix"
$ws.Cells.Item(142, 11).Value = 1.8
$ws.Cells.Item(142, 12).Value = 3.6
$ws.Cells.Item(142, 13).Value = 4.333
$ws.Cells.Item(142, 14).Value = 1.727
$ws.Cells.Item(142, 15).Value = 3.6
$ws.Cells.Item(142, 16).Value = 4.75
$ws.Cells.Item(142, 17).Value = -0.75
$ws.Cells.Item(142, 18).Value = 1.99
$ws.Cells.Item(142, 19).Value = 1.91
$ws.Cells.Item(142, 20).Value = 2.75
$ws.Cells.Item(142, 21).Value = 2
$ws.Cells.Item(142, 22).Value = 1.85
$ws.Cells.Item(142, 23).Value = 0
$ws.Cells.Item(142, 24).Value = 0
$ws.Cells.Item(142, 25).Value = 0
$ws.Cells.Item(142, 26).Value = 0
$ws.Cells.Item(142, 27).Value = 0

# Row 143
$ws.Range("A134").Copy()
$ws.Range("A143").PasteSpecial(-4122)
$ws.Range("E134").Copy()
$ws.Range("E143").PasteSpecial(-4122)
$ws.Cells.Item(143, 1).Value = 141
$ws.Cells.Item(143, 2).Value = 7126794
$ws.Cells.Item(143, 3).Value = "Australia ALeague"
$ws.Cells.Item(143, 4).Value = "Australia ALeague"
$ws.Cells.Item(143, 5).Value = 45388.23958333334
$ws.Cells.Item(143, 6).Value = "Melbourne Victory"
$ws.Cells.Item(143, 7).Value = "Melbourne City"
$ws.Cells.Item(143, 11).Value = 1.833
$ws.Cells.Item(143, 12).Value = 3.5
$ws.Cells.Item(143, 13).Value = 4.5
$ws.Cells.Item(143, 14).Value = 2.15
$ws.Cells.Item(143, 15).Value = 3.75
$ws.Cells.Item(143, 16).Value = 3
$ws.Cells.Item(143, 17).Value = -0.25
$ws.Cells.Item(143, 18).Value = 1.95
$ws.Cells.Item(143, 19).Value = 1.95
$ws.Cells.Item(143, 20).Value = 2.75
$ws.Cells.Item(143, 21).Value = 1.8
$ws.Cells.Item(143, 22).Value = 2.05
$ws.Cells.Item(143, 23).Value = 0
$ws.Cells.Item(143, 24).Value = 0
$ws.Cells.Item(143, 25).Value = 0
$ws.Cells.Item(143, 26).Value = 0
$ws.Cells.Item(143, 27).Value = 0

# Row 144
$ws.Range("A134").Copy()
$ws.Range("A144").PasteSpecial(-4122)
$ws.Range("E134").Copy()
$ws.Range("E144").PasteSpecial(-4122)
$ws.Cells.Item(144, 1).Value = 142
$ws.Cells.Item(144, 2).Value = 7127403
$ws.Cells.Item(144, 3).Value = "Australia ALeague"
$ws.Cells.Item(144, 4).Value = "Australia ALeague"
$ws.Cells.Item(144, 5).Value = 45389.08333333334
$ws.Cells.Item(144, 6).Value = "Newcastle Jets"
$ws.Cells.Item(144, 7).Value = "Sydney FC"
$ws.Cells.Item(144, 11).Value = 3.6
$ws.Cells.Item(144, 12).Value = 3.6
$ws.Cells.Item(144, 13).Value = 1.952
$ws.Cells.Item(144, 14).Value = 4
$ws.Cells.Item(144, 15).Value = 3.75
$ws.Cells.Item(144, 16).Value = 1.85
$ws.Cells.Item(144, 17).Value = 0.5
$ws.Cells.Item(144, 18).Value = 1.98
$ws.Cells.Item(144, 19).Value = 1.92
$ws.Cells.Item(144, 20).Value = 3
$ws.Cells.Item(144, 21).Value = 1.875
$ws.Cells.Item(144, 22).Value = 1.975
$ws.Cells.Item(144, 23).Value = 0
$ws.Cells.Item(144, 24).Value = 0
$ws.Cells.Item(144, 25).Value = 0
$ws.Cells.Item(144, 26).Value = 0
$ws.Cells.Item(144, 27).Value = 0

# Row 145
$ws.Range("A134").Copy()
$ws.Range("A145").PasteSpecial(-4122)
$ws.Range("E134").Copy()
$ws.Range("E145").PasteSpecial(-4122)
$ws.Cells.Item(145, 1).Value = 143
$ws.Cells.Item(145, 2).Value = 7127402
$ws.Cells.Item(145, 3).Value = "Australia ALeague"
$ws.Cells.Item(145, 4).Value = "Australia ALeague"
$ws.Cells.Item(145, 5).Value = 45389.16666666666
$ws.Cells.Item(145, 6).Value = "Perth Glory"
$ws.Cells.Item(145, 7).Value = "Adelaide United"
$ws.Cells.Item(145, 11).Value = 2.25
$ws.Cells.Item(145, 12).Value = 3.5
$ws.Cells.Item(145, 13).Value = 3
$ws.Cells.Item(145, 14).Value = 2.6
$ws.Cells.Item(145, 15).Value = 3.5
$ws.Cells.Item(145, 16).Value = 2.6
$ws.Cells.Item(145, 17).Value = 0
$ws.Cells.Item(145, 18).Value = 1.95
$ws.Cells.Item(145, 19).Value = 1.95
$ws.Cells.Item(145, 20).Value = 3.25
$ws.Cells.Item(145, 21).Value = 1.925
$ws.Cells.Item(145, 22).Value = 1.925
$ws.Cells.Item(145, 23).Value = 0
$ws.Cells.Item(145, 24).Value = 0
$ws.Cells.Item(145, 25).Value = 0
$ws.Cells.Item(145, 26).Value = 0
$ws.Cells.Item(145, 27).Value = 0

Write-Output "done"